$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A. This shifts the existing
# Name / Email / Department / Delegate Comments / Predicted Category
# columns one place to the right (B:F) and keeps their data/formatting.
$ws.Columns("A").Insert()

# New header for the inserted Employee_ID column (bold, matching the
# other header cells in row 1).
$ws.Range("A1").Value = "Employee_ID"
$ws.Range("A1").Font.Bold = $true

# New data value for the Employee_ID column on the existing data row.
# Force text formatting so the leading zero is preserved ("001" and
# not the number 1), then copy the plain formatting from a neighboring
# data cell so no stray number format is left behind on the cell.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "001"
$ws.Range("B2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Add a new (sparse) record row for Eve Adams, only populating the
# Name and Predicted Category columns.
$ws.Range("B3").Value = "Eve Adams"
$ws.Range("F3").Value = "Add"

# Column widths for the new layout (columns C:F keep the widths they
# had before the insert; A and B take on new auto-fit-style widths).
# (Values chosen land on the closest representable column width.)
$ws.Columns("A").ColumnWidth = 12.0
$ws.Columns("B").ColumnWidth = 9.666666666666668
